$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.721958666666667
$ws.Range("H2").Value = 5.165876000000001
$ws.Range("I2").Value = 0.07789986924239836
$ws.Range("J2").Value = 0.07789986924239835
$ws.Range("M2").Value = 154.942487
$ws.Range("N2").Value = 464.827461
$ws.Range("O2").Value = 0.982851703624775
$ws.Range("P2").Value = 0.9828517036247751
$ws.Range("Q2").Value = 266.8045583245374
$ws.Range("R2").Value = 2401.241024920836
$ws.Range("S2").Value = 0.07656401919703844
$ws.Range("T2").Value = 0.07656401919703844
$ws.Range("G3").Value = 1.721958666666667
$ws.Range("H3").Value = 5.165876000000001
$ws.Range("I3").Value = 0.07789986924239836
$ws.Range("J3").Value = 0.07789986924239835
$ws.Range("O3").Value = 0.003358739549735124
$ws.Range("P3").Value = 0.003358739549735124
$ws.Range("Q3").Value = 0.9117621903582223
$ws.Range("R3").Value = 8.205859713224001
$ws.Range("S3").Value = 0.0002616453717436381
$ws.Range("T3").Value = 0.0002616453717436381
$ws.Range("G4").Value = 1.721958666666667
$ws.Range("H4").Value = 5.165876000000001
$ws.Range("I4").Value = 0.07789986924239836
$ws.Range("J4").Value = 0.07789986924239835
$ws.Range("M4").Value = 1.771368666666667
$ws.Range("N4").Value = 5.314106000000001
$ws.Range("O4").Value = 0.01123638032078883
$ws.Range("P4").Value = 0.01123638032078884
$ws.Range("Q4").Value = 3.050223627428445
$ws.Range("R4").Value = 27.45201264685601
$ws.Range("S4").Value = 0.0008753125577473084
$ws.Range("T4").Value = 0.0008753125577473082
$ws.Range("G5").Value = 1.721958666666667
$ws.Range("H5").Value = 5.165876000000001
$ws.Range("I5").Value = 0.07789986924239836
$ws.Range("J5").Value = 0.07789986924239835
$ws.Range("M5").Value = 0.4024976666666666
$ws.Range("N5").Value = 1.207493
$ws.Range("O5").Value = 0.002553176504700935
$ws.Range("P5").Value = 0.002553176504700936
$ws.Range("Q5").Value = 0.6930843454297778
$ws.Range("R5").Value = 6.237759108868
$ws.Range("S5").Value = 0.0001988921158689666
$ws.Range("T5").Value = 0.0001988921158689666
$ws.Range("I6").Value = 0.3978297504389287
$ws.Range("J6").Value = 0.3978297504389286
$ws.Range("M6").Value = 154.942487
$ws.Range("N6").Value = 464.827461
$ws.Range("O6").Value = 0.982851703624775
$ws.Range("P6").Value = 0.9828517036247751
$ws.Range("Q6").Value = 1362.554159416345
$ws.Range("R6").Value = 12262.9874347471
$ws.Range("S6").Value = 0.3910076479715201
$ws.Range("T6").Value = 0.3910076479715201
$ws.Range("I7").Value = 0.3978297504389287
$ws.Range("J7").Value = 0.3978297504389286
$ws.Range("O7").Value = 0.003358739549735124
$ws.Range("P7").Value = 0.003358739549735124
$ws.Range("S7").Value = 0.001336206516860484
$ws.Range("T7").Value = 0.001336206516860484
$ws.Range("I8").Value = 0.3978297504389287
$ws.Range("J8").Value = 0.3978297504389286
$ws.Range("M8").Value = 1.771368666666667
$ws.Range("N8").Value = 5.314106000000001
$ws.Range("O8").Value = 0.01123638032078883
$ws.Range("P8").Value = 0.01123638032078884
$ws.Range("Q8").Value = 15.57730091570333
$ws.Range("R8").Value = 140.19570824133
$ws.Range("S8").Value = 0.004470166378856311
$ws.Range("T8").Value = 0.004470166378856311
$ws.Range("I9").Value = 0.3978297504389287
$ws.Range("J9").Value = 0.3978297504389286
$ws.Range("M9").Value = 0.4024976666666666
$ws.Range("N9").Value = 1.207493
$ws.Range("O9").Value = 0.002553176504700935
$ws.Range("P9").Value = 0.002553176504700936
$ws.Range("Q9").Value = 3.539538318318333
$ws.Range("R9").Value = 31.855844864865
$ws.Range("S9").Value = 0.001015729571691709
$ws.Range("T9").Value = 0.001015729571691709
$ws.Range("G10").Value = 5.953778333333333
$ws.Range("H10").Value = 17.861335
$ws.Range("I10").Value = 0.269343604258924
$ws.Range("J10").Value = 0.269343604258924
$ws.Range("M10").Value = 154.942487
$ws.Range("N10").Value = 464.827461
$ws.Range("O10").Value = 0.982851703624775
$ws.Range("P10").Value = 0.9828517036247751
$ws.Range("Q10").Value = 922.4932220133817
$ws.Range("R10").Value = 8302.438998120435
$ws.Range("S10").Value = 0.2647248203063206
$ws.Range("T10").Value = 0.2647248203063207
$ws.Range("G11").Value = 5.953778333333333
$ws.Range("H11").Value = 17.861335
$ws.Range("I11").Value = 0.269343604258924
$ws.Range("J11").Value = 0.269343604258924
$ws.Range("O11").Value = 0.003358739549735124
$ws.Range("P11").Value = 0.003358739549735124
$ws.Range("Q11").Value = 3.152474028087778
$ws.Range("R11").Value = 28.37226625279
$ws.Range("S11").Value = 0.0009046550160926537
$ws.Range("T11").Value = 0.0009046550160926537
$ws.Range("G12").Value = 5.953778333333333
$ws.Range("H12").Value = 17.861335
$ws.Range("I12").Value = 0.269343604258924
$ws.Range("J12").Value = 0.269343604258924
$ws.Range("M12").Value = 1.771368666666667
$ws.Range("N12").Value = 5.314106000000001
$ws.Range("O12").Value = 0.01123638032078883
$ws.Range("P12").Value = 0.01123638032078884
$ws.Range("Q12").Value = 10.54633638794556
$ws.Range("R12").Value = 94.91702749151001
$ws.Range("S12").Value = 0.003026447174425309
$ws.Range("T12").Value = 0.003026447174425309
$ws.Range("G13").Value = 5.953778333333333
$ws.Range("H13").Value = 17.861335
$ws.Range("I13").Value = 0.269343604258924
$ws.Range("J13").Value = 0.269343604258924
$ws.Range("M13").Value = 0.4024976666666666
$ws.Range("N13").Value = 1.207493
$ws.Range("O13").Value = 0.002553176504700935
$ws.Range("P13").Value = 0.002553176504700936
$ws.Range("Q13").Value = 2.396381887017222
$ws.Range("R13").Value = 21.567436983155
$ws.Range("S13").Value = 0.0006876817620853514
$ws.Range("T13").Value = 0.0006876817620853515
$ws.Range("G14").Value = 5.635097666666667
$ws.Range("H14").Value = 16.905293
$ws.Range("I14").Value = 0.2549267760597491
$ws.Range("J14").Value = 0.254926776059749
$ws.Range("M14").Value = 154.942487
$ws.Range("N14").Value = 464.827461
$ws.Range("O14").Value = 0.982851703624775
$ws.Range("P14").Value = 0.9828517036247751
$ws.Range("Q14").Value = 873.1160469612304
$ws.Range("R14").Value = 7858.044422651073
$ws.Range("S14").Value = 0.2505552161498959
$ws.Range("T14").Value = 0.2505552161498958
$ws.Range("G15").Value = 5.635097666666667
$ws.Range("H15").Value = 16.905293
$ws.Range("I15").Value = 0.2549267760597491
$ws.Range("J15").Value = 0.254926776059749
$ws.Range("O15").Value = 0.003358739549735124
$ws.Range("P15").Value = 0.003358739549735124
$ws.Range("Q15").Value = 2.983735376986889
$ws.Range("R15").Value = 26.853618392882
$ws.Range("S15").Value = 0.0008562326450383483
$ws.Range("T15").Value = 0.0008562326450383482
$ws.Range("G16").Value = 5.635097666666667
$ws.Range("H16").Value = 16.905293
$ws.Range("I16").Value = 0.2549267760597491
$ws.Range("J16").Value = 0.254926776059749
$ws.Range("M16").Value = 1.771368666666667
$ws.Range("N16").Value = 5.314106000000001
$ws.Range("O16").Value = 0.01123638032078883
$ws.Range("P16").Value = 0.01123638032078884
$ws.Range("Q16").Value = 9.981835440339779
$ws.Range("R16").Value = 89.83651896305801
$ws.Range("S16").Value = 0.002864454209759906
$ws.Range("T16").Value = 0.002864454209759906
$ws.Range("G17").Value = 5.635097666666667
$ws.Range("H17").Value = 16.905293
$ws.Range("I17").Value = 0.2549267760597491
$ws.Range("J17").Value = 0.254926776059749
$ws.Range("M17").Value = 0.4024976666666666
$ws.Range("N17").Value = 1.207493
$ws.Range("O17").Value = 0.002553176504700935
$ws.Range("P17").Value = 0.002553176504700936
$ws.Range("Q17").Value = 2.268113662272111
$ws.Range("R17").Value = 20.413022960449
$ws.Range("S17").Value = 0.0006508730550549082
$ws.Range("T17").Value = 0.0006508730550549082
